# "arrumando os champs certo" — fix the ASSISTS (F) column so it stores
# real numbers instead of text-as-number, and correct a batch of
# mis-labelled CHAMPION (H) cells to "Rakan".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F column (ASSISTS): convert text-number cells F2:F41 to real numbers ---
$assists = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0;
    7  = 2;  8  = 2;  9  = 2;  10 = 2;  11 = 2;  12 = 2;  13 = 2;
    14 = 5;  15 = 5;  16 = 5;  17 = 5;
    18 = 8;
    19 = 10; 20 = 10; 21 = 10; 22 = 10;
    23 = 11; 24 = 11;
    25 = 12; 26 = 12; 27 = 12;
    28 = 13; 29 = 13;
    30 = 14;
    31 = 15; 32 = 15; 33 = 15; 34 = 15; 35 = 15;
    36 = 17;
    37 = 18; 38 = 18;
    39 = 19; 40 = 19; 41 = 19
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# --- H column (CHAMPION): relabel a batch of rows to "Rakan" ---
$champions = @(5, 11, 17, 23, 25, 29, 35, 36, 41)

foreach ($row in $champions) {
    $ws.Cells.Item($row, 8).Value = "Rakan"
}
